# Bottom-up Plan Microwave - fix mistake in row 6 (Step 4)
# The "X" marker in column E (Light) of step 4 was placed there by mistake;
# clear the cell's content while leaving its formatting/style untouched.

$wb = $excel.ActiveWorkbook

if ($wb.Worksheets.Count -ge 1) {
    $ws = $wb.Worksheets.Item(1)
} else {
    $ws = $wb.ActiveSheet
}

$ws.Range("E6").ClearContents()
